$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Remove "generosi.me, " from the "(see: nathanrlapierre.com, generosi.me,
# cnfsr.com, nanbendo.com)" list -- the site was taken down.
#
# A plain Range.Delete() across that run collapses every run-of-identical-
# -formatting "w:r" sibling that follows the deletion point into a single
# run (losing their distinct w:rsidR attributes/run boundaries). To keep the
# surrounding runs intact exactly as they were, we temporarily drop a zero-
# length bookmark at the start of each of those following runs before
# deleting -- a bookmark anchored inside a run "splits" it and blocks the
# delete-time coalescing -- then remove those scratch bookmarks again (which
# does not re-merge anything, it just drops the bookmarkStart/End pair).
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$needle = "generosi.me, "
$delStart = $full.IndexOf($needle)
$delEnd = $delStart + $needle.Length

# Starting text of each run that must stay a distinct run after the delete.
$followingRunStarts = @(
    "cnfsr.com, nanbendo.com)",
    "; ",
    "have some",
    " experience with webservers and apache "
)

$scratchBookmarks = New-Object System.Collections.ArrayList
$n = 0
foreach ($runStart in $followingRunStarts) {
    $pos = $full.IndexOf($runStart, $delEnd)
    $n = $n + 1
    $name = "zzRunGuard" + $n
    $guardRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $guardRange)
    [void]$scratchBookmarks.Add($name)
}

# Delete the "generosi.me, " text itself.
$deleteRange = $d.Range($delStart, $delEnd)
$deleteRange.Delete()

# Drop the scratch guard bookmarks now that the delete is done.
foreach ($name in $scratchBookmarks) {
    $d.Bookmarks.Item($name).Delete()
}

# Word tracks the location of the last edit with a "_GoBack" bookmark (and a
# document can only have one). Planting it at the point of our deletion
# mirrors that and automatically replaces the stale one that used to sit
# after the GPA "3.92".
$goBackRange = $d.Range($delStart, $delStart)
$d.Bookmarks.Add("_GoBack", $goBackRange)
